$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamp values for column Z (rows 2-48), per the commit diff.
$timestamps = @(
    "2025-10-19T23:55:19.370719",
    "2025-10-19T23:55:19.371721",
    "2025-10-19T23:55:19.371721",
    "2025-10-19T23:55:19.371721",
    "2025-10-19T23:55:19.371721",
    "2025-10-19T23:55:19.372724",
    "2025-10-19T23:55:19.372724",
    "2025-10-19T23:55:19.372724",
    "2025-10-19T23:55:19.372724",
    "2025-10-19T23:55:19.372724",
    "2025-10-19T23:55:19.373724",
    "2025-10-19T23:55:19.373724",
    "2025-10-19T23:55:19.373724",
    "2025-10-19T23:55:19.373724",
    "2025-10-19T23:55:19.414454",
    "2025-10-19T23:55:19.414454",
    "2025-10-19T23:55:19.414454",
    "2025-10-19T23:55:19.414454",
    "2025-10-19T23:55:19.415438",
    "2025-10-19T23:55:19.415438",
    "2025-10-19T23:55:19.415438",
    "2025-10-19T23:55:19.415438",
    "2025-10-19T23:55:19.415438",
    "2025-10-19T23:55:19.416438",
    "2025-10-19T23:55:19.503201",
    "2025-10-19T23:55:19.504166",
    "2025-10-19T23:55:19.504166",
    "2025-10-19T23:55:19.505170",
    "2025-10-19T23:55:19.505170",
    "2025-10-19T23:55:19.506167",
    "2025-10-19T23:55:19.506167",
    "2025-10-19T23:55:19.506167",
    "2025-10-19T23:55:19.507166",
    "2025-10-19T23:55:19.507166",
    "2025-10-19T23:55:19.507166",
    "2025-10-19T23:55:19.508167",
    "2025-10-19T23:55:19.510361",
    "2025-10-19T23:55:19.510892",
    "2025-10-19T23:55:19.510892",
    "2025-10-19T23:55:19.511892",
    "2025-10-19T23:55:19.511892",
    "2025-10-19T23:55:19.511892",
    "2025-10-19T23:55:19.512899",
    "2025-10-19T23:55:19.513894",
    "2025-10-19T23:55:19.513894",
    "2025-10-19T23:55:19.513894",
    "2025-10-19T23:55:19.513894"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 26).Value = $timestamps[$i]
}
